$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet: the "Procedure:" block (rows 9-11) is replaced by a new
# "Notes:" block (rows 9-10), and the original "Procedure:" block is moved
# down to rows 12-14.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("About")

# Former bold "Procedure:" header at A9 becomes a plain "Notes:" header.
$ws1.Range("A9").Value = "Notes:"
$ws1.Range("A9").Font.Bold = $false

# Former procedure description at A10 becomes the notes description.
$ws1.Range("A10").Value = "This variable sets the fraction of O&M costs in the electricity sector that is labor."

# Former A11 text ("which are best representative...") moves away from here.
$ws1.Range("A11").ClearContents()

# Re-insert the "Procedure:" block (still bold) further down, at rows 12-14.
$ws1.Range("A12").Value = "Procedure:"
$ws1.Range("A12").Font.Bold = $true
$ws1.Range("A13").Value = "Divide ""Fixed Labor Costs"" by ""Total Fixed O&M Costs.""  Using values for supercritical boilers,"
$ws1.Range("A14").Value = "which are best representative of the type of plant that would be built going forward."

[void]$ws1.Activate()
[void]$ws1.Range("A11").Select()

# ---------------------------------------------------------------------------
# "FoOMCtiL" sheet: the header label gains a "(dimensionless)" suffix and
# wraps onto two lines.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("FoOMCtiL")
$ws2.Range("B1").Value = "Frac of O&M Costs (dimensionless)"
$ws2.Range("B1").WrapText = $true
$ws2.Rows.Item(1).RowHeight = 30

[void]$ws2.Activate()
[void]$ws2.Range("B1").Select()

[void]$ws1.Activate()
